$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update status text for all anomaly rows (shared string used by C2:C5)
$ws.Range("C2:C5").Value = "Reparada e encerrada."

# Update "Data da última modificação" (column D) from 2015-06-09 to 2015-06-10 for rows 2-5
$newDate = Get-Date -Year 2015 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("D2:D5").Value = $newDate
